$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Score")

# Fill in the missing "O" column (homework/extra-score) values for a handful
# of students, which ripples through the existing P (avg), Q (total) and
# R (grade lookup) formulas already in the sheet.
$ws.Range("O14").Value = 41
$ws.Range("O19").Value = 28
$ws.Range("O24").Value = 40
$ws.Range("O28").Value = 25
$ws.Range("O30").Value = 38
$ws.Range("O31").Value = 40

# Column widths: columns P and Q (16 & 17) now share the same width (5
# characters) instead of P being narrower than Q. 4.14 chars of
# Excel's ColumnWidth serializes to the stored width="5" seen in the diff.
$ws.Columns.Item(16).ColumnWidth = 4.14
$ws.Columns.Item(17).ColumnWidth = 4.14

# View state: scrolled down a bit further and a different active selection.
$ws.Range("O15").Select()

# Slightly smaller print scale.
$ws.PageSetup.Zoom = 86
